$wb = $excel.ActiveWorkbook

# New row (60) appends 2025-04-30 to each "Silver/wafer price" sheet, carrying
# forward the last known price (same value as row 59, the 2025-04-29 row).
# The USD_CNY sheet (row 59 is the last date sheet untouched) is NOT updated.
$updates = @(
    @{ Sheet = "N-Dense";                   Price = "38" },
    @{ Sheet = "N-Type";                    Price = "37.3" },
    @{ Sheet = "N-type Wafer";              Price = "1.03" },
    @{ Sheet = "Cell Topcon 183mm";         Price = "0.275" },
    @{ Sheet = "Module Topcon 183mm";       Price = "0.09" },
    @{ Sheet = "Silver Rear_side";          Price = "5,375" },
    @{ Sheet = "Silver Busbar front-side";  Price = "8,047" },
    @{ Sheet = "Silver finger front-side";  Price = "8,097" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    # Leading apostrophe forces text entry (matches existing text-typed cells
    # in columns A/B) instead of Excel's automatic date/number coercion.
    $ws.Range("A60").Value = "'2025-04-30"
    $ws.Range("B60").Value = "'" + $u.Price
}
